$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# Log Week 13 entry for the new WR player
$ws.Cells.Item(10, 1).Value = "J.Mickens"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0

$ws.Range("J11").Select()
